# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-10-07 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-08 Sunday", 2) | Out-Null

# Update the division-problem table. Addressing cells directly by
# (row, column) avoids any ambiguity from values that coincide with other
# cells' old/new text (e.g. "10÷7=1, 3" is both an old value in row 13
# and the new value written into row 5).
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "49÷8=6, 1" },
    @{ Row = 1;  Col = 2; Text = "27÷7=3, 6" },
    @{ Row = 1;  Col = 3; Text = "79÷6=13, 1" },
    @{ Row = 1;  Col = 4; Text = "51÷3=17, 0" },
    @{ Row = 1;  Col = 5; Text = "89÷5=17, 4" },

    @{ Row = 5;  Col = 1; Text = "51÷4=12, 3" },
    @{ Row = 5;  Col = 2; Text = "10÷7=1, 3" },
    @{ Row = 5;  Col = 3; Text = "19÷3=6, 1" },
    @{ Row = 5;  Col = 4; Text = "57÷3=19, 0" },
    @{ Row = 5;  Col = 5; Text = "55÷2=27, 1" },

    @{ Row = 9;  Col = 1; Text = "50÷7=7, 1" },
    @{ Row = 9;  Col = 2; Text = "48÷2=24, 0" },
    @{ Row = 9;  Col = 3; Text = "21÷2=10, 1" },
    @{ Row = 9;  Col = 4; Text = "85÷2=42, 1" },
    @{ Row = 9;  Col = 5; Text = "19÷8=2, 3" },

    @{ Row = 13; Col = 1; Text = "16÷6=2, 4" },
    @{ Row = 13; Col = 2; Text = "75÷2=37, 1" },
    @{ Row = 13; Col = 3; Text = "76÷2=38, 0" },
    @{ Row = 13; Col = 4; Text = "64÷4=16, 0" },
    @{ Row = 13; Col = 5; Text = "88÷9=9, 7" },

    @{ Row = 17; Col = 1; Text = "75÷7=10, 5" },
    @{ Row = 17; Col = 2; Text = "55÷6=9, 1" },
    @{ Row = 17; Col = 3; Text = "20÷2=10, 0" },
    @{ Row = 17; Col = 4; Text = "31÷6=5, 1" },
    @{ Row = 17; Col = 5; Text = "59÷7=8, 3" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
